$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 used to describe "agents/personnesMorales" -> now "agents/collectivites"
$ws.Range("A4").Value = "agents/collectivites"

# Both the "collectivites" (row 4) and "personnesPhysiques" (row 5) referentiels now
# mention that there is one file per agent, and their rows grew taller to fit the text.
$ws.Range("G4").Value = "Référentiel produit automatiquement à partir du RI_013 du SIA. Pas de différence notable avec le contenu du référentiel SIA. Un fichier par agent. 40 entités alignées (owl:sameAs) avec autant d'entités du référentiel des producteurs."
$ws.Range("G5").Value = "Référentiel produit automatiquement à partir du RI_012 du SIA. Pas de différence notable avec le contenu du référentiel SIA. Un fichier par agent. 102 entités alignées (owl:sameAs) avec autant d'entités du référentiel des producteurs."

$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 60

# Move the active selection to G4
$ws.Range("G4").Select()
